# Commit: "update actions - relative direction"
# Row 2 of Sheet1 stores per-feature attribution scores for this trajectory.
# This edit overwrites the scores that changed with their new values
# (re-run of the attribution computation), leaving all other cells,
# formatting and the rest of the sheet untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.08608604857666213
$ws.Range("D2").Value = 0.2707430686247467
$ws.Range("E2").Value = 0.004001273688284372
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.05711919423029652
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2424185176332642
$ws.Range("N2").Value = 0.01439279748655964
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1039938171523975
$ws.Range("V2").Value = 0.01705956103150915
$ws.Range("W2").Value = -0.01078891192103048
$ws.Range("Y2").Value = -0
$ws.Range("Z2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = [double]"-3.470842752942556e-11"
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.006287610896074329
$ws.Range("AF2").Value = [double]"-1.992104022987901e-12"
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.04586640419228727
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.04249581251515201
$ws.Range("AO2").Value = 0.07080012785383533
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AS2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1675629633249888
$ws.Range("AW2").Value = 0.09601824582853803
$ws.Range("AX2").Value = -0.01667735686670889
$ws.Range("AY2").Value = -0
$ws.Range("BB2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.007563276064319164
$ws.Range("BF2").Value = 0.1049887927948767
$ws.Range("BG2").Value = 0.004555816567386153
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.01855919159624354
$ws.Range("BO2").Value = -0.03870817360056667
$ws.Range("BP2").Value = -0.08014909847657475
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.06240362818182101
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.02238961253711967
$ws.Range("BY2").Value = -0.02819281016782701
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03017968464745394
$ws.Range("CG2").Value = -0.04373135439866592
$ws.Range("CH2").Value = 0.01555492970330826
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01683645637809813
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.02927517608387242
$ws.Range("CQ2").Value = 0.05732951286099017
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04485442205455682
$ws.Range("CY2").Value = -0.04610213201090916
$ws.Range("CZ2").Value = 0.007004487323011817
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03663010099907456
$ws.Range("DH2").Value = 0.01496390616943506
$ws.Range("DI2").Value = 0.03101452545121341
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.02091649284860365
$ws.Range("DQ2").Value = 0.04285491153331428
$ws.Range("DR2").Value = -0.03731491465691677
$ws.Range("DS2").Value = -0
$ws.Range("DU2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05320440195189079
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.002892739116739693
$ws.Range("EA2").Value = -0.03046001104498006
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04214883444383082
$ws.Range("EH2").Value = 0
$ws.Range("EI2").Value = 0.09694243965812084
$ws.Range("EJ2").Value = -0.02535434639816732
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.0523474947198154
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.0505857322256974
$ws.Range("ES2").Value = 0.02099614862398629
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.03977438091833668
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.03758484801388312
$ws.Range("FB2").Value = 0.01722156883334725
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.006986232390577869
$ws.Range("FJ2").Value = -0.0160849485541347
$ws.Range("FK2").Value = 0.02278434274117715
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.009903314644151768
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = 0.01161490472161882
$ws.Range("FT2").Value = -0.00335253712351679
$ws.Range("FU2").Value = 0
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03276389726031452
$ws.Range("GB2").Value = 0.02249852629173317
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
